$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force the written value to stay a TEXT cell (matches the
# original inline-string cells) instead of Excel's automatic number
# inference, while leaving the cell's style index untouched (reset back
# to "Normal" after the write so no new style entry is introduced).
function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "56.407.28"
Set-TextCell "E2" "  +3.23%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.972.77"
Set-TextCell "E3" "  +2.73%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.03%  "

# Row 5 - BNB
Set-TextCell "D5" "501.96"
Set-TextCell "E5" "  +5.18%  "

# Row 6 - Solana
Set-TextCell "D6" "134.73"
Set-TextCell "E6" "  +5.88%  "

# Row 8 - XRP
Set-TextCell "E8" "  +5.21%  "

# Row 9 - Toncoin
Set-TextCell "D9" "7.44"
Set-TextCell "E9" "  +10.96%  "

# Row 10 - Dogecoin
Set-TextCell "E10" "  +9.77%  "

# Row 11 - Cardano
Set-TextCell "E11" "  +4.10%  "

# Row 12 - TRON
Set-TextCell "E12" "  +3.52%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextCell "D13" "3.482.74"
Set-TextCell "E13" "  +2.81%  "

# Row 14 - Avalanche
Set-TextCell "D14" "25.41"
Set-TextCell "E14" "  +10.84%  "

# Row 15 - ShibaInu
Set-TextCell "E15" "  +11.16%  "

# Row 16 - WrappedBTC
Set-TextCell "D16" "56.407.34"
Set-TextCell "E16" "  +3.37%  "

# Row 17 - WrappedEther
Set-TextCell "D17" "2.976.08"
Set-TextCell "E17" "  +2.95%  "

# Row 18 - Polkadot
Set-TextCell "D18" "5.75"
Set-TextCell "E18" "  +8.89%  "

# Row 19 - Chainlink
Set-TextCell "D19" "12.31"
Set-TextCell "E19" "  +5.99%  "

# Row 20 - Uniswap
Set-TextCell "D20" "7.72"
Set-TextCell "E20" "  +8.68%  "

# Row 21 - BitcoinCash
Set-TextCell "D21" "321.33"
Set-TextCell "E21" "  +4.25%  "

# Row 22 - Dai
Set-TextCell "E22" "  -0.35%  "

# Row 23 - Polygon
Set-TextCell "E23" "  +4.25%  "

# Row 24 - Litecoin
Set-TextCell "D24" "61.65"
Set-TextCell "E24" "  +3.25%  "

# Row 25 - Binance-PegBSC-USD
Set-TextCell "D25" "1.00"
Set-TextCell "E25" "  +1.21%  "

# Row 26 - Kaspa
Set-TextCell "E26" "  +4.58%  "

# Row 27 - PEPE
Set-TextCell "D27" "0.0₃0887"
Set-TextCell "E27" "  +8.13%  "

# Row 28 - RenderToken
Set-TextCell "E28" "  +1.18%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextCell "D29" "6.74"
Set-TextCell "E29" "  +8.24%  "

# Row 30 - Fetch.AI
Set-TextCell "D30" "1.18"
Set-TextCell "E30" "  +3.18%  "

# Row 31 - PancakeSwap
Set-TextCell "E31" "  +7.50%  "

# Row 32 - EthereumClassic
Set-TextCell "D32" "20.31"
Set-TextCell "E32" "  +6.11%  "

# Row 33 - Monero
Set-TextCell "D33" "158.27"
Set-TextCell "E33" "  +15.95%  "

# Row 34 - NEARProtocol
Set-TextCell "E34" "  +4.01%  "

# Row 35 - now ImmutableX (was Aptos)
Set-TextCell "B35" "ImmutableX"
Set-TextCell "C35" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D35" "1.25"
Set-TextCell "E35" "  +2.71%  "

# Row 36 - now Aptos (was ImmutableX)
Set-TextCell "B36" "Aptos"
Set-TextCell "C36" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D36" "5.50"
Set-TextCell "E36" "  +0.94%  "

# Row 37 - Hedera
Set-TextCell "D37" "0.0671"
Set-TextCell "E37" "  +7.92%  "

# Row 38 - EnergySwap
Set-TextCell "D38" "23.05"
Set-TextCell "E38" "  +0.41%  "

# Row 39 - RenzoRestakedETH
Set-TextCell "D39" "3.007.06"
Set-TextCell "E39" "  +2.99%  "

# Row 40 - FirstDigitalUSD
Set-TextCell "E40" "  +0.00%  "

# Row 41 - OKB
Set-TextCell "D41" "36.21"
Set-TextCell "E41" "  +1.29%  "

# Row 42 - Mantle
Set-TextCell "D42" "0.640"
Set-TextCell "E42" "  +5.80%  "

# Row 43 - Maker
Set-TextCell "D43" "2.238.65"
Set-TextCell "E43" "  +8.88%  "

# Row 45 - ONDO
Set-TextCell "D45" "0.977"
Set-TextCell "E45" "  +0.68%  "

# Row 46 - Filecoin
Set-TextCell "D46" "3.55"
Set-TextCell "E46" "  +3.29%  "

# Row 47 - dogwifhat
Set-TextCell "E47" "  +18.23%  "

# Row 48 - VeChain
Set-TextCell "E48" "  +9.08%  "

# Row 49 - Cosmos
Set-TextCell "E49" "  +7.21%  "

# Row 50 - InjectiveProtocol
Set-TextCell "D50" "18.85"
Set-TextCell "E50" "  +4.92%  "

# Row 51 - Stellar
Set-TextCell "E51" "  +8.69%  "
